$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.2860218016679856
$ws.Range("D2").Value = 0.7775373733460067

$ws.Range("C3").Value = 0.766301375340139
$ws.Range("D3").Value = 0.4516416712719507

$ws.Range("C4").Value = -0.1847189028023646
$ws.Range("D4").Value = 0.8551417512401858

$ws.Range("C5").Value = -0.001178267201751197
$ws.Range("D5").Value = 0.9990704981126601

$ws.Range("C6").Value = 0.8822070828919951
$ws.Range("D6").Value = 0.3872027250654795

$ws.Range("C7").Value = 0.1494013801071472
$ws.Range("D7").Value = 0.8825980610306894

$ws.Range("C8").Value = 0.2334529381475392
$ws.Range("D8").Value = 0.8175682481166846

$ws.Range("C9").Value = -0.9191922406697574
$ws.Range("D9").Value = 0.3679640350462561

$ws.Range("C10").Value = -0.5589123310282735
$ws.Range("D10").Value = 0.5818667697948272

$ws.Range("C11").Value = 0.20962497102554
$ws.Range("D11").Value = 0.8358894500962069
